# ---------------------------------------------------------------------------
# PlayerPerformance_3434.xlsx edit
#
# 1. Insert a new "Player Info" sheet before "ODI Batting" with player bio.
# 2. Rename "MATCH_CARD_LINK" -> "MATCH_CODE" on "ODI Batting" / "ODI Bowling"
#    and replace the scraped howstat.com URL with the bare numeric match code.
# 3. Drop the now-pointless empty INNING_NUMBER cells on "ODI Batting".
# 4. Append a new "ODI Batting Extra" sheet with extra per-match batting
#    detail (batting position, boundary counts, % of total runs, MoM).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# --- 1. "Player Info" sheet, inserted ahead of "ODI Batting" ---------------

$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
$playerInfo.Range("A1:D1").Style = $battingSheet.Range("A1").Style

$playerInfo.Cells.Item(2, 1).NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value2 = "3434"
$playerInfo.Range("B2").Value = "Stuart Christopher John Broad"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium Fast"

# --- 2 & 3. "ODI Batting": header rename, URL -> code, drop empty cells ----

$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"

$lastRow = $battingSheet.Range("A1").End(4).Row
for ($r = 2; $r -le $lastRow; $r++) {

    $inning = $battingSheet.Cells.Item($r, 2).Value()
    if ([string]::IsNullOrEmpty($inning)) {
        $battingSheet.Cells.Item($r, 2).ClearContents()
    }

    $link = $battingSheet.Cells.Item($r, 4).Value()
    if ($link -match "MatchCode=(\d+)") {
        $code = $matches[1]
        $cell = $battingSheet.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value2 = $code
    }
}

# --- 2. "ODI Bowling": header rename, URL -> code ---------------------------

$bowlingSheet.Cells.Item(1, 2).Value = "MATCH_CODE"

$lastRow = $bowlingSheet.Range("A1").End(4).Row
for ($r = 2; $r -le $lastRow; $r++) {

    $link = $bowlingSheet.Cells.Item($r, 2).Value()
    if ($link -match "MatchCode=(\d+)") {
        $code = $matches[1]
        $cell = $bowlingSheet.Cells.Item($r, 2)
        $cell.NumberFormat = "@"
        $cell.Value2 = $code
    }
}

# --- 4. "ODI Batting Extra" sheet, appended after "ODI Bowling" ------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$extra.Name = "ODI Batting Extra"

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"
$extra.Range("A1:F1").Style = $battingSheet.Range("A1").Style

$extraRows = @(
    @("3519", "",   "",  "",  "",       "NO"),
    @("3602", 9,    "0", "0", "0.41%",  "NO"),
    @("3606", "",   "",  "",  "",       "NO"),
    @("3608", 9,    "1", "0", "3.30%",  "NO"),
    @("3622", 10,   "",  "",  "",       "NO"),
    @("3625", 9,    "3", "0", "17.18%", "NO"),
    @("3629", 9,    "",  "",  "",       "NO"),
    @("3727", 10,   "0", "0", "",       "NO"),
    @("3735", 9,    "",  "",  "",       "NO"),
    @("3738", "",   "",  "",  "",       "NO"),
    @("3744", 9,    "0", "0", "1.49%",  "NO"),
    @("3746", "",   "",  "",  "",       "NO"),
    @("3749", 9,    "0", "0", "",       "NO"),
    @("3756", 9,    "0", "0", "3.25%",  "NO"),
    @("3761", "",   "",  "",  "",       "NO"),
    @("3769", 9,    "",  "",  "",       "NO"),
    @("3780", 10,   "0", "1", "3.46%",  "NO"),
    @("3785", 10,   "",  "",  "",       "NO"),
    @("3890", 10,   "1", "0", "2.29%",  "NO"),
    @("3891", 10,   "0", "1", "5.51%",  "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $codeCell = $extra.Cells.Item($r, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value2 = $row[0]

    if (-not [string]::IsNullOrEmpty($row[1])) {
        $extra.Cells.Item($r, 2).Value = $row[1]
    }
    if (-not [string]::IsNullOrEmpty($row[2])) {
        $extra.Cells.Item($r, 3).Value = $row[2]
    }
    if (-not [string]::IsNullOrEmpty($row[3])) {
        $extra.Cells.Item($r, 4).Value = $row[3]
    }
    if (-not [string]::IsNullOrEmpty($row[4])) {
        $extra.Cells.Item($r, 5).Value = $row[4]
    }
    $extra.Cells.Item($r, 6).Value = $row[5]

    $r = $r + 1
}
